# Singapore Premier League workbook update (28-06-2024 19:47 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Rows 20 and 21 were re-sorted (same match date, 2023-07-02):
#    swap every column except A (the running index) between the
#    two rows so the Hougang/Balestier match now sits on row 20
#    and the Geylang/Tampines match on row 21.
# ---------------------------------------------------------------
$ws.Range("B20").Value = 6228597
$ws.Range("E20").Value = "Hougang United FC"
$ws.Range("F20").Value = "Balestier Khalsa FC"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = "A"
$ws.Range("L20").Value = 2.5
$ws.Range("M20").Value = 3.6
$ws.Range("N20").Value = 2.25
$ws.Range("O20").Value = 2.6
$ws.Range("P20").Value = 3.75
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 0.25
$ws.Range("S20").Value = 1.825
$ws.Range("T20").Value = 2.025
$ws.Range("U20").Value = 4
$ws.Range("V20").Value = 1.95
$ws.Range("W20").Value = 1.9
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = -1
$ws.Range("Z20").Value = 1.2
$ws.Range("AA20").Value = -1
$ws.Range("AB20").Value = 1.025
$ws.Range("AC20").Value = 0
$ws.Range("AD20").Value = 0

$ws.Range("B21").Value = 6228598
$ws.Range("E21").Value = "Geylang International"
$ws.Range("F21").Value = "Tampines Rovers FC"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = "D"
$ws.Range("L21").Value = 3.6
$ws.Range("M21").Value = 4.2
$ws.Range("N21").Value = 1.666
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 4.5
$ws.Range("Q21").Value = 1.55
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1.85
$ws.Range("T21").Value = 2
$ws.Range("U21").Value = 3.25
$ws.Range("V21").Value = 2.025
$ws.Range("W21").Value = 1.825
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = 3.5
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = 0.8500000000000001
$ws.Range("AB21").Value = -1
$ws.Range("AC21").Value = -1
$ws.Range("AD21").Value = 0.825

# ---------------------------------------------------------------
# 2) Append the new, not-yet-played fixture as row 70
#    (Tanjong Pagar United vs DPMM FC, 2024-06-29).
# ---------------------------------------------------------------
$ws.Range("A70").Value = 68
$ws.Range("A69").Copy()
$ws.Range("A70").PasteSpecial(-4122)

# B70 holds the raw match id, but it was written out as text (it
# carries the same shared-string typing the source feed used for
# this particular row) rather than as a plain number.
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = "8089721"
$ws.Range("B70").Style = "Normal"

$ws.Range("C70").Value = "Singapore Premier League"

$ws.Range("D70").Value = 45472.29166666666
$ws.Range("D69").Copy()
$ws.Range("D70").PasteSpecial(-4122)

$ws.Range("E70").Value = "Tanjong Pagar United"
$ws.Range("F70").Value = "DPMM FC"

# No result yet (G70:K70 intentionally left blank) - only the
# opening/closing odds columns are populated.
$ws.Range("L70").Value = 3.75
$ws.Range("M70").Value = 4.1
$ws.Range("N70").Value = 1.65
$ws.Range("O70").Value = 3.7
$ws.Range("P70").Value = 3.9
$ws.Range("Q70").Value = 1.7
$ws.Range("R70").Value = 0.75
$ws.Range("S70").Value = 1.9
$ws.Range("T70").Value = 1.95
$ws.Range("U70").Value = 3.5
$ws.Range("V70").Value = 1.875
$ws.Range("W70").Value = 1.975
$ws.Range("X70").Value = 0
$ws.Range("Y70").Value = 0
$ws.Range("Z70").Value = 0
